$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Range("E7").Value = "Jennifer Espich, Christina Frank, Daniel Lackmann, Granit Gecaj, Masood Ahmed, Markus Schmidtner"
$ws.Range("E8").Value = "Jennifer Espich, Christina Frank, Daniel Lackmann, Granit Gecaj, Masood Ahmed, Markus Schmidtner"
